$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for Case_2_188 (380 kV case)
# Rows 2-25 correspond to timesteps 0-23; only columns with non-zero
# loading percentages changed (B, C, E, F, G, I, K, L, M, N).

$ws.Range("B2").Value = 16.82891876849569
$ws.Range("C2").Value = 4.953342742468764
$ws.Range("E2").Value = 9.513429722921821
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.748348860011498
$ws.Range("I2").Value = 34.07006900356462
$ws.Range("K2").Value = 14.9732497098677
$ws.Range("L2").Value = 10.4909847853254
$ws.Range("M2").Value = 16.52295696420623
$ws.Range("N2").Value = 24.12936905743975

$ws.Range("B3").Value = 16.70224838573854
$ws.Range("C3").Value = 4.800533114539136
$ws.Range("E3").Value = 9.526393896011234
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.751059598505672
$ws.Range("I3").Value = 34.12895254745397
$ws.Range("K3").Value = 14.88964861149583
$ws.Range("L3").Value = 10.50128802694281
$ws.Range("M3").Value = 16.51709261662836
$ws.Range("N3").Value = 24.18789383821448

$ws.Range("B4").Value = 16.62812510105525
$ws.Range("C4").Value = 4.705857100406258
$ws.Range("E4").Value = 9.53517072756082
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.75281128651042
$ws.Range("I4").Value = 34.17046250735546
$ws.Range("K4").Value = 14.84160306294024
$ws.Range("L4").Value = 10.50906765749515
$ws.Range("M4").Value = 16.51630827292374
$ws.Range("N4").Value = 24.22578865046102

$ws.Range("B5").Value = 16.59886685984464
$ws.Range("C5").Value = 4.667130321994076
$ws.Range("E5").Value = 9.53895314581117
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.753547136263625
$ws.Range("I5").Value = 34.18872282330026
$ws.Range("K5").Value = 14.822866809608
$ws.Range("L5").Value = 10.51260369065644
$ws.Range("M5").Value = 16.51669869714317
$ws.Range("N5").Value = 24.24172477486781

$ws.Range("B6").Value = 16.59406659913082
$ws.Range("C6").Value = 4.660693053778921
$ws.Range("E6").Value = 9.539593654321907
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.753670655860041
$ws.Range("I6").Value = 34.19183608966169
$ws.Range("K6").Value = 14.81980702081239
$ws.Range("L6").Value = 10.51321294552508
$ws.Range("M6").Value = 16.51680645974533
$ws.Range("N6").Value = 24.24440079503946

$ws.Range("B7").Value = 16.62772664029974
$ws.Range("C7").Value = 4.705335313611728
$ws.Range("E7").Value = 9.535220904807248
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.752821121170449
$ws.Range("I7").Value = 34.17070333059504
$ws.Range("K7").Value = 14.84134694649144
$ws.Range("L7").Value = 10.50911386433178
$ws.Range("M7").Value = 16.51631066116576
$ws.Range("N7").Value = 24.22600157076883

$ws.Range("B8").Value = 16.78450262442965
$ws.Range("C8").Value = 4.900870142618877
$ws.Range("E8").Value = 9.51773049055668
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.74926545185391
$ws.Range("I8").Value = 34.08925944553728
$ws.Range("K8").Value = 14.94375232131751
$ws.Range("L8").Value = 10.49423588747573
$ws.Range("M8").Value = 16.52035159435954
$ws.Range("N8").Value = 24.14914197026571

$ws.Range("B9").Value = 17.11956823436754
$ws.Range("C9").Value = 5.274784701256792
$ws.Range("E9").Value = 9.48989490446862
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.742981951316242
$ws.Range("I9").Value = 33.97212268631245
$ws.Range("K9").Value = 15.1698596817624
$ws.Range("L9").Value = 10.47657758474091
$ws.Range("M9").Value = 16.55052989657061
$ws.Range("N9").Value = 24.01394117575819

$ws.Range("B10").Value = 17.38067032939437
$ws.Range("C10").Value = 5.540293237598821
$ws.Range("E10").Value = 9.473360370894307
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.738780830745832
$ws.Range("I10").Value = 33.9121247447244
$ws.Range("K10").Value = 15.3502740941061
$ws.Range("L10").Value = 10.47060337090745
$ws.Range("M10").Value = 16.58612707575501
$ws.Range("N10").Value = 23.9240227648727

$ws.Range("B11").Value = 17.50226338162067
$ws.Range("C11").Value = 5.658459123869014
$ws.Range("E11").Value = 9.466683648067731
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.736958811082551
$ws.Range("I11").Value = 33.89050907186437
$ws.Range("K11").Value = 15.43519223238012
$ws.Range("L11").Value = 10.46939941205967
$ws.Range("M11").Value = 16.60519938121108
$ws.Range("N11").Value = 23.88515062402723

$ws.Range("B12").Value = 17.54867513134705
$ws.Range("C12").Value = 5.702784124431385
$ws.Range("E12").Value = 9.464276417198123
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.736281593643355
$ws.Range("I12").Value = 33.88314146919159
$ws.Range("K12").Value = 15.46773464385863
$ws.Range("L12").Value = 10.4691605306824
$ws.Range("M12").Value = 16.61283182591211
$ws.Range("N12").Value = 23.8707223063772

$ws.Range("B13").Value = 17.53866381227504
$ws.Range("C13").Value = 5.693257384164841
$ws.Range("E13").Value = 9.464789477935485
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.73642687893701
$ws.Range("I13").Value = 33.88469182399812
$ws.Range("K13").Value = 15.46070927967217
$ws.Range("L13").Value = 10.46920233654244
$ws.Range("M13").Value = 16.6111698640275
$ws.Range("N13").Value = 23.87381673634834

$ws.Range("B14").Value = 17.50607456577564
$ws.Range("C14").Value = 5.662114468643705
$ws.Range("E14").Value = 9.466483179032368
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.736902841026197
$ws.Range("I14").Value = 33.88988653823495
$ws.Range("K14").Value = 15.43786192264147
$ws.Range("L14").Value = 10.46937541345123
$ws.Range("M14").Value = 16.60581910954117
$ws.Range("N14").Value = 23.88395775401908

$ws.Range("B15").Value = 17.48615938532749
$ws.Range("C15").Value = 5.642982316430691
$ws.Range("E15").Value = 9.467536378608116
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.737196038865616
$ws.Range("I15").Value = 33.89317498642014
$ws.Range("K15").Value = 15.42391678137533
$ws.Range("L15").Value = 10.4695096716506
$ws.Range("M15").Value = 16.60259491598154
$ws.Range("N15").Value = 23.89020739576907

$ws.Range("B16").Value = 17.37277738339142
$ws.Range("C16").Value = 5.532514172291806
$ws.Range("E16").Value = 9.473813678365969
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.738901690820448
$ws.Range("I16").Value = 33.91365174858169
$ws.Range("K16").Value = 15.34477977611982
$ws.Range("L16").Value = 10.47071245957793
$ws.Range("M16").Value = 16.58493826954458
$ws.Range("N16").Value = 23.92660399968478

$ws.Range("B17").Value = 17.30391576002643
$ws.Range("C17").Value = 5.464041639566227
$ws.Range("E17").Value = 9.477880720171864
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.739970821950586
$ws.Range("I17").Value = 33.92766878980083
$ws.Range("K17").Value = 15.29694433780358
$ws.Range("L17").Value = 10.47183760298377
$ws.Range("M17").Value = 16.57484125524374
$ws.Range("N17").Value = 23.94945226263015

$ws.Range("B18").Value = 17.26457577758344
$ws.Range("C18").Value = 5.424414320291131
$ws.Range("E18").Value = 9.480299527094349
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.740594147694444
$ws.Range("I18").Value = 33.93626536508273
$ws.Range("K18").Value = 15.26970006680871
$ws.Range("L18").Value = 10.47262727931237
$ws.Range("M18").Value = 16.56930495194987
$ws.Range("N18").Value = 23.96278528795704

$ws.Range("B19").Value = 17.251303007765
$ws.Range("C19").Value = 5.410956802081759
$ws.Range("E19").Value = 9.481132168768196
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.740806638157605
$ws.Range("I19").Value = 33.93926774775585
$ws.Range("K19").Value = 15.26052260320193
$ws.Range("L19").Value = 10.47291914706909
$ws.Range("M19").Value = 16.56747714989876
$ws.Range("N19").Value = 23.96733249460025

$ws.Range("B20").Value = 17.31121879511204
$ws.Range("C20").Value = 5.471356229150119
$ws.Range("E20").Value = 9.477439546174358
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.739856143307106
$ws.Range("I20").Value = 33.92612133883253
$ws.Range("K20").Value = 15.30200878786372
$ws.Range("L20").Value = 10.471703083306
$ws.Range("M20").Value = 16.57588805468301
$ws.Range("N20").Value = 23.94700023000744

$ws.Range("B21").Value = 17.51563715889503
$ws.Range("C21").Value = 5.671273688668465
$ws.Range("E21").Value = 9.465982414576372
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.736762694162671
$ws.Range("I21").Value = 33.88833852069405
$ws.Range("K21").Value = 15.44456246739
$ws.Range("L21").Value = 10.46931869187955
$ws.Range("M21").Value = 16.60737965585225
$ws.Range("N21").Value = 23.88097117809977

$ws.Range("B22").Value = 17.65135823888489
$ws.Range("C22").Value = 5.799454411921276
$ws.Range("E22").Value = 9.459200186363558
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.734815183806364
$ws.Range("I22").Value = 33.86841240601918
$ws.Range("K22").Value = 15.53996661519234
$ws.Range("L22").Value = 10.4690250302503
$ws.Range("M22").Value = 16.63035005271267
$ws.Range("N22").Value = 23.83951751307373

$ws.Range("B23").Value = 17.57873982887333
$ws.Range("C23").Value = 5.731282514612515
$ws.Range("E23").Value = 9.462755553638051
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.73584783679148
$ws.Range("I23").Value = 33.87861076662543
$ws.Range("K23").Value = 15.48885089971681
$ws.Range("L23").Value = 10.46906627661261
$ws.Range("M23").Value = 16.61787309069246
$ws.Range("N23").Value = 23.86148671421395

$ws.Range("B24").Value = 17.30791631292837
$ws.Range("C24").Value = 5.468050114343635
$ws.Range("E24").Value = 9.477638749830938
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.739907962515586
$ws.Range("I24").Value = 33.92681926560886
$ws.Range("K24").Value = 15.2997183477693
$ws.Range("L24").Value = 10.47176345474405
$ws.Range("M24").Value = 16.57541395967325
$ws.Range("N24").Value = 23.94810817941435

$ws.Range("B25").Value = 17.02615898538344
$ws.Range("C25").Value = 5.17502047041009
$ws.Range("E25").Value = 9.496735692408999
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.74460851963587
$ws.Range("I25").Value = 33.99924153421123
$ws.Range("K25").Value = 15.10609594592057
$ws.Range("L25").Value = 10.48012360255818
$ws.Range("M25").Value = 16.53999667662035
$ws.Range("N25").Value = 24.04886003815091
